# Update NONE icons on cards: card image URLs moved from the raw
# githubusercontent "master/images" path to the GitHub Pages
# "choice_cards/cards" path. Replace the common URL prefix across every
# cell on the sheet (the card_url_1 .. card_url_6 columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldPrefix = "https://raw.githubusercontent.com/bryanparthum/farmland_conservation/master/images/"
$newPrefix = "https://bryanparthum.github.io/farmland_conservation/choice_cards/cards/"

# LookAt:=xlPart (2) so the match is a substring replace, not a whole-cell
# match (cells contain the prefix plus a per-card file name suffix).
$ws.Cells.Replace($oldPrefix, $newPrefix, 2)
